# LeetCodeStats.xlsx - add "Power of Two" entry to the List sheet and
# refresh the view/selection state to match (recursion approach was
# fantastic!).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "List"
$ws2 = $wb.Worksheets.Item(2)   # "Statistic"

# --- New "Power of Two" row (row 59) -------------------------------------
# Add the hyperlink first (anchored on the still-empty cell), then set the
# display text + reapply the built-in "Hyperlink" style so the cell reuses
# the workbook's existing hyperlink style (xf index 5) instead of minting a
# near-duplicate one.
$ws1.Hyperlinks.Add($ws1.Range("B59"), "https://leetcode.com/problems/power-of-two/", "", "", "https://leetcode.com/problems/power-of-two/")
$ws1.Range("B59").Value = "Power of Two"
$ws1.Range("B59").Style = "Hyperlink"

$ws1.Cells.Item(59, 3).Value = 1
$ws1.Cells.Item(59, 4).Value = 1
$ws1.Cells.Item(59, 9).Value = "https://leetcode.com/problems/power-of-two/submissions/1093210199/"

# --- Placeholder rows 60-77 (column A index only, like the existing tail) -
for ($i = 60; $i -le 77; $i++) {
    $ws1.Cells.Item($i, 1).Value = $i - 2
}

# --- View / selection state -----------------------------------------------
# Before: "List" tab selected, selection C67. After: "Statistic" tab
# selected, "List" selection moved to B68, "Statistic" selection moved to
# E17.
$ws1.Range("B68").Select()
$ws2.Activate()
$ws2.Range("E17").Select()
